$wb = $excel.ActiveWorkbook

# --- Sheet "Hoja1": update the conversion summary text in A1 ---
$wsHoja1 = $wb.Worksheets.Item("Hoja1")
$wsHoja1.Range("A1").Value = "Conversión del día 💰`n✅ Dólar paralelo: 68`n`nBinance`n✅ 1000 Bs = 5.13 = 20333.44 pesos`n✅ 20333.44 pesos = 5.11 = 953.16 Bs`n`nPromedio competencia`n✅ Tasa pesos: 20`n✅ Tasa Bs: 20`n✅ % Ganancia: 20%"

# --- Sheet "tasas": update rate values ---
$wsTasas = $wb.Worksheets.Item("tasas")
$wsTasas.Range("N10").Value = 194.999
$wsTasas.Range("O10").Value = 3965
$wsTasas.Range("N12").Value = 3982.8
$wsTasas.Range("O12").Value = 186.7
